# Update ConfigMotos report: rename/reshuffle columns, drop the old
# "Data de Cadastro" column I (table shrinks from A1:I2 to A1:H2), and
# fill in a full data row describing the registered motorcycle.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ConfigMotos")

# Remove the now-unused last column (I) entirely, shifting nothing in
# because it's the rightmost column of the used range.
$ws.Columns.Item(9).Delete()

# Header row (row 1) - new column layout: nome, marca, cor, nome_dono,
# observacoes, valor_compra, status, Data de Cadastro
$ws.Range("A1").Value = "nome"
$ws.Range("B1").Value = "marca"
$ws.Range("C1").Value = "cor"
$ws.Range("D1").Value = "nome_dono"
$ws.Range("E1").Value = "observacoes"
$ws.Range("F1").Value = "valor_compra"
$ws.Range("G1").Value = "status"
$ws.Range("H1").Value = "Data de Cadastro"

# Data row (row 2)
$ws.Range("A2").Value = "Yamaha MT-07"
$ws.Range("B2").Value = "Yamaha"
$ws.Range("C2").Value = "preto"
$ws.Range("D2").Value = "Matheus"
$ws.Range("E2").Value = "teste"
$ws.Range("F2").Value = 15000
$ws.Range("G2").Value = "Ativo"
$ws.Range("H2").Value = "2024-05-17 16:22:23"

# Shrink the autofilter / filter-database range to match the new extent.
# (Toggle the stale A1:I2 filter off first - calling AutoFilter() on a
# range while a filter is already active would just remove it instead of
# resizing it.)
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}
$ws.Range("A1:H2").AutoFilter()

# Keep the workbook-level _FilterDatabase defined name in sync with the
# new (narrower) autofilter range.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*FilterDatabase*") {
        $n.RefersTo = "='ConfigMotos'!`$A`$1:`$H`$2"
    }
}
